# Adds summary statistics rows to the Winner-determination results sheet:
#   J12  = average of the k column (|S*|/n), bold
#   A14:B14 = "Average of SW(S*)/SW(OPT)" / =AVERAGE(N2:N11)
#   A15:B15 = "Average of SC(S*)/SC(OPT)" / =AVERAGE(Z2:Z11)
#   A16:B16 = "Worst of SW(S*)/SW(OPT)"   / =MIN(N2:N11)
#   A17:B17 = "Worst of SC(S*)/SC(OPT)"   / =MAX(Z2:Z11)
# B14:B17 are bold, size 12, vertically centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New labeled summary rows (A14:B17) -----------------------------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format B14 (bold, size 12, vertical-center), then clone the format onto
# B15:B17 via copy/paste-special so every cell lands on the same style index
# instead of each re-deriving its own (avoids orphan style entries).
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- J12: average of the k column (column J), bold -------------------------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection (matches the saved cursor position in the workbook) ---------
$ws.Range("A14:B17").Select()
